# Generate Report for Handback
# - Status for the 5c9ecb06... row flips from "Ready for handoff" to
#   "Handback transform failed" on both the zh-cn and de-de sheets.
# - The Error Detail column (P) gets a message explaining the handback
#   file-name mismatch for each locale.
# - Column P is widened to fit the new error text.

$wb = $excel.ActiveWorkbook

$ws_zhcn = $wb.Worksheets.Item("zh-cn")
$ws_dede = $wb.Worksheets.Item("de-de")

# Status column (C) on row 7 - the 5c9ecb06... entry
$ws_zhcn.Range("C7").Value = "Handback transform failed"
$ws_dede.Range("C7").Value = "Handback transform failed"

# Error Detail column (P) on row 7
$ws_zhcn.Range("P7").Value = "Handback file name: uefg3xel.yr1 is different with handoff file name: 5c9ecb06-4feb-4f20-ba40-d8056ccc1ba4.dcd4b84cbe6f02c9700e12153ee33a861d5d6c54.zh-cn."
$ws_dede.Range("P7").Value = "Handback file name: uefg3xel.yr1 is different with handoff file name: 5c9ecb06-4feb-4f20-ba40-d8056ccc1ba4.dcd4b84cbe6f02c9700e12153ee33a861d5d6c54.de-de."

# Widen column P so the new error text is readable.
# NOTE: Excel's ColumnWidth property is in "characters" and gets converted
# to the internal width unit with font-metric padding, so 39.17 here is
# what round-trips to a stored width of exactly 40 (matching the other
# width="40" columns on this sheet, e.g. column A/G/I/J).
$ws_zhcn.Columns.Item(16).ColumnWidth = 39.17
$ws_dede.Columns.Item(16).ColumnWidth = 39.17
